$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 23.77523808933262
$ws.Range("C2").Value = 16.50886831778488
$ws.Range("D2").Value = 6.065269279119201
$ws.Range("E2").Value = 7.354145249903395
$ws.Range("G2").Value = 3.75874884933868
$ws.Range("I2").Value = 50.36420617299364
$ws.Range("L2").Value = 11.29222126842757

$ws.Range("B3").Value = 23.6145350902627
$ws.Range("C3").Value = 15.88631421548643
$ws.Range("D3").Value = 5.962698343820747
$ws.Range("E3").Value = 7.335720397160289
$ws.Range("G3").Value = 3.765156196712545
$ws.Range("I3").Value = 49.04475038904852
$ws.Range("L3").Value = 11.27627308599382

$ws.Range("B4").Value = 23.52852340840468
$ws.Range("C4").Value = 15.49813573381704
$ws.Range("D4").Value = 5.901054266776604
$ws.Range("E4").Value = 7.324232198230381
$ws.Range("G4").Value = 3.769277306055131
$ws.Range("I4").Value = 48.21913970308682
$ws.Range("L4").Value = 11.26919269501724

$ws.Range("B5").Value = 23.49668179905144
$ws.Range("C5").Value = 15.33873381940032
$ws.Range("D5").Value = 5.876301385331165
$ws.Range("E5").Value = 7.319505883819788
$ws.Range("G5").Value = 3.771003979864485
$ws.Range("I5").Value = 47.87915159728737
$ws.Range("L5").Value = 11.2669887387207

$ws.Range("B6").Value = 23.49158891767241
$ws.Range("C6").Value = 15.31220030618667
$ws.Range("D6").Value = 5.872214307179093
$ws.Range("E6").Value = 7.318718357669058
$ws.Range("G6").Value = 3.771293556792553
$ws.Range("I6").Value = 47.82249310304139
$ws.Range("L6").Value = 11.26666389986712

$ws.Range("B7").Value = 23.52808096106527
$ws.Range("C7").Value = 15.49599052683381
$ws.Range("D7").Value = 5.900718910805269
$ws.Range("E7").Value = 7.324168639750024
$ws.Range("G7").Value = 3.769300400751649
$ws.Range("I7").Value = 48.21456839179891
$ws.Range("L7").Value = 11.2691602139771

$ws.Range("B8").Value = 23.71721746215704
$ws.Range("C8").Value = 16.29559909644534
$ws.Range("D8").Value = 6.02964510495941
$ws.Range("E8").Value = 7.347828305510065
$ws.Range("G8").Value = 3.760919461643765
$ws.Range("I8").Value = 49.91266566457823
$ws.Range("L8").Value = 11.28615929626019

$ws.Range("B9").Value = 24.18709167193156
$ws.Range("C9").Value = 17.80579195070857
$ws.Range("D9").Value = 6.291634356827049
$ws.Range("E9").Value = 7.392867790128599
$ws.Range("G9").Value = 3.745955500938267
$ws.Range("I9").Value = 53.10598084163969
$ws.Range("L9").Value = 11.34101732494315

$ws.Range("B10").Value = 24.59017498007402
$ws.Range("C10").Value = 18.8673625526625
$ws.Range("D10").Value = 6.487789658531641
$ws.Range("E10").Value = 7.425180005837763
$ws.Range("G10").Value = 3.735840646515481
$ws.Range("I10").Value = 55.35164364851484
$ws.Range("L10").Value = 11.39442249353165

$ws.Range("B11").Value = 24.78545509151069
$ws.Range("C11").Value = 19.33767882967233
$ws.Range("D11").Value = 6.57742378670462
$ws.Range("E11").Value = 7.439721509711927
$ws.Range("G11").Value = 3.731426210584606
$ws.Range("I11").Value = 56.34825886873703
$ws.Range("L11").Value = 11.4215481139618

$ws.Range("B12").Value = 24.86105305158911
$ws.Range("C12").Value = 19.51381018726279
$ws.Range("D12").Value = 6.611391453260744
$ws.Range("E12").Value = 7.445206196435469
$ws.Range("G12").Value = 3.729781144919806
$ws.Range("I12").Value = 56.72182777822054
$ws.Range("L12").Value = 11.43222487932738

$ws.Range("B13").Value = 24.84469937222048
$ws.Range("C13").Value = 19.47596690228212
$ws.Range("D13").Value = 6.604075296749327
$ws.Range("E13").Value = 7.444025930908182
$ws.Range("G13").Value = 3.730134261417188
$ws.Range("I13").Value = 56.64154685253616
$ws.Range("L13").Value = 11.42990748087298

$ws.Range("B14").Value = 24.79164190114597
$ws.Range("C14").Value = 19.35220950723212
$ws.Range("D14").Value = 6.580218005583414
$ws.Range("E14").Value = 7.440173171869399
$ws.Range("G14").Value = 3.731290338874986
$ws.Range("I14").Value = 56.37907061588203
$ws.Range("L14").Value = 11.42241839813753

$ws.Range("B15").Value = 24.75935552949789
$ws.Range("C15").Value = 19.27614407056051
$ws.Range("D15").Value = 6.565607053981322
$ws.Range("E15").Value = 7.437810421660601
$ws.Range("G15").Value = 3.73200192339322
$ws.Range("I15").Value = 56.21779112269155
$ws.Range("L15").Value = 11.41788377076649

$ws.Range("B16").Value = 24.5776473820963
$ws.Range("C16").Value = 18.83635899998912
$ws.Range("D16").Value = 6.481937131819723
$ws.Range("E16").Value = 7.424226581794197
$ws.Range("G16").Value = 3.736132876115791
$ws.Range("I16").Value = 55.28598944329334
$ws.Range("L16").Value = 11.39270656216826

$ws.Range("B17").Value = 24.46918319786797
$ws.Range("C17").Value = 18.5632211534727
$ws.Range("D17").Value = 6.430687751510665
$ws.Range("E17").Value = 7.415853551533254
$ws.Range("G17").Value = 3.738714739639791
$ws.Range("I17").Value = 54.70779153067393
$ws.Range("L17").Value = 11.37798488225111

$ws.Range("B18").Value = 24.4079231415301
$ws.Range("C18").Value = 18.40494168194963
$ws.Range("D18").Value = 6.40125020150821
$ws.Range("E18").Value = 7.411022735185433
$ws.Range("G18").Value = 3.740217364810503
$ws.Range("I18").Value = 54.37289207779168
$ws.Range("L18").Value = 11.36978402788222

$ws.Range("B19").Value = 24.38737674507372
$ws.Range("C19").Value = 18.35115376908259
$ws.Range("D19").Value = 6.391290991062826
$ws.Range("E19").Value = 7.409384526173097
$ws.Range("G19").Value = 3.740729160283417
$ws.Range("I19").Value = 54.25910746998674
$ws.Range("L19").Value = 11.36705322055048

$ws.Range("B20").Value = 24.48061334079644
$ws.Range("C20").Value = 18.59242023553322
$ws.Range("D20").Value = 6.436139462920135
$ws.Range("E20").Value = 7.416746408355049
$ws.Range("G20").Value = 3.738438075809372
$ws.Range("I20").Value = 54.76958519123366
$ws.Range("L20").Value = 11.37952444606321

$ws.Range("B21").Value = 24.80718192612443
$ws.Range("C21").Value = 19.38861458289812
$ws.Range("D21").Value = 6.587225041765675
$ws.Range("E21").Value = 7.441305407714542
$ws.Range("G21").Value = 3.730950051666559
$ws.Range("I21").Value = 56.45627193384301
$ws.Range("L21").Value = 11.42460715535649

$ws.Range("B22").Value = 25.03019616748784
$ws.Range("C22").Value = 19.89744911335831
$ws.Range("D22").Value = 6.686100612790227
$ws.Range("E22").Value = 7.457229059735688
$ws.Range("G22").Value = 3.726211020840331
$ws.Range("I22").Value = 57.53621295781244
$ws.Range("L22").Value = 11.45642988109103

$ws.Range("B23").Value = 24.91031459508922
$ws.Range("C23").Value = 19.62697525907284
$ws.Range("D23").Value = 6.63332722264586
$ws.Range("E23").Value = 7.448741672313199
$ws.Range("G23").Value = 3.728726259132987
$ws.Range("I23").Value = 56.96195189974718
$ws.Range("L23").Value = 11.43923055702044

$ws.Range("B24").Value = 24.47544235008341
$ws.Range("C24").Value = 18.57922320849811
$ws.Range("D24").Value = 6.433674659466821
$ws.Range("E24").Value = 7.41634280095053
$ws.Range("G24").Value = 3.738563098575501
$ws.Range("I24").Value = 54.74165600145238
$ws.Range("L24").Value = 11.37882759042377

$ws.Range("B25").Value = 24.04961133173084
$ws.Range("C25").Value = 17.40484352031753
$ws.Range("D25").Value = 6.219976708203681
$ws.Range("E25").Value = 7.380823668383869
$ws.Range("G25").Value = 3.749847977274703
$ws.Range("I25").Value = 52.2585918099022
$ws.Range("L25").Value = 11.32387263133237
